# Aggiornamento dati fino al 9 agosto 2021
# Appends rows 329:343 (dates 44403..44417) to the daily series.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 329
$dates = 44403,44404,44405,44406,44407,44408,44409,44410,44411,44412,44413,44414,44415,44416,44417
$bVals  = 0,0,0,0,0,0,0,0,0,0,0,1,0,1,0
$cVals  = 0,0,0,0,0,0,0,0,0,0,0,1,1,2,2
$dVals  = 0,0,0,0,0,0,0,0,0,0,0,28.87669650591972,28.87669650591972,57.75339301183945,57.75339301183945

for ($i = 0; $i -lt $dates.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $dates[$i]
    $ws.Cells.Item($r, 2).Value = $bVals[$i]
    $ws.Cells.Item($r, 3).Value = $cVals[$i]
    $ws.Cells.Item($r, 4).Value = $dVals[$i]
}

$lastRow = $startRow + $dates.Length - 1

# Mirror the formatting (date number format, bold font, border, centered
# alignment) of the preceding data row onto the newly appended rows, same
# as dragging the fill handle down would do.
$ws.Range("A328:D328").Copy() | Out-Null
$ws.Range("A" + $startRow + ":D" + $lastRow).PasteSpecial(-4122) | Out-Null
